# Auto-generated edit script: applies the South Korea K League 1 update
# (team-name swap Gangwon FC <-> Incheon Utd, plus row-level data rotations)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("F4").Value = 'Incheon Utd'
# Row 6
$ws.Range("F6").Value = 'Gangwon FC'
# Row 10
$ws.Range("B10").Value = 6149376
$ws.Range("E10").Value = 'Jeonbuk Motors'
$ws.Range("F10").Value = 'Daegu FC'
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 1.75
$ws.Range("M10").Value = 3.5
$ws.Range("N10").Value = 4.2
$ws.Range("O10").Value = 2.05
$ws.Range("P10").Value = 3.3
$ws.Range("Q10").Value = 3.4
$ws.Range("R10").Value = -0.25
$ws.Range("S10").Value = 1.825
$ws.Range("T10").Value = 2.025
$ws.Range("U10").Value = 2.25
$ws.Range("V10").Value = 1.85
$ws.Range("W10").Value = 2
$ws.Range("X10").Value = 1.05
$ws.Range("AA10").Value = 0.825
$ws.Range("AC10").Value = -1
$ws.Range("AD10").Value = 1
# Row 11
$ws.Range("B11").Value = 6149854
$ws.Range("E11").Value = 'Gwangju FC'
$ws.Range("F11").Value = 'Suwon Bluewings'
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("L11").Value = 1.8
$ws.Range("M11").Value = 3.4
$ws.Range("N11").Value = 3.8
$ws.Range("O11").Value = 1.85
$ws.Range("P11").Value = 3.6
$ws.Range("Q11").Value = 4.2
$ws.Range("R11").Value = -0.5
$ws.Range("S11").Value = 1.9
$ws.Range("T11").Value = 1.95
$ws.Range("U11").Value = 2.5
$ws.Range("V11").Value = 2.025
$ws.Range("W11").Value = 1.825
$ws.Range("X11").Value = 0.8500000000000001
$ws.Range("AA11").Value = 0.8999999999999999
$ws.Range("AC11").Value = 1.025
$ws.Range("AD11").Value = -1
# Row 12
$ws.Range("B12").Value = 6149855
$ws.Range("F12").Value = 'FC Seoul'
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = 1
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 'D'
$ws.Range("L12").Value = 2.55
$ws.Range("M12").Value = 3.1
$ws.Range("N12").Value = 2.625
$ws.Range("O12").Value = 3
$ws.Range("P12").Value = 3.1
$ws.Range("Q12").Value = 2.25
$ws.Range("R12").Value = 0.25
$ws.Range("S12").Value = 1.825
$ws.Range("T12").Value = 2.025
$ws.Range("U12").Value = 2.25
$ws.Range("V12").Value = 1.825
$ws.Range("W12").Value = 2.025
$ws.Range("Y12").Value = 2.1
$ws.Range("Z12").Value = -1
$ws.Range("AA12").Value = 0.4125
$ws.Range("AB12").Value = -0.5
$ws.Range("AC12").Value = -0.5
$ws.Range("AD12").Value = 0.5125
# Row 13
$ws.Range("B13").Value = 6149856
$ws.Range("F13").Value = 'Daejeon Hana Citizen'
$ws.Range("H13").Value = 2
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1
$ws.Range("K13").Value = 'A'
$ws.Range("L13").Value = 2.3
$ws.Range("M13").Value = 3.2
$ws.Range("N13").Value = 2.875
$ws.Range("O13").Value = 2.375
$ws.Range("P13").Value = 3.3
$ws.Range("Q13").Value = 2.7
$ws.Range("R13").Value = 0
$ws.Range("S13").Value = 1.775
$ws.Range("T13").Value = 2.1
$ws.Range("U13").Value = 2.5
$ws.Range("V13").Value = 1.95
$ws.Range("W13").Value = 1.9
$ws.Range("Y13").Value = -1
$ws.Range("Z13").Value = 1.7
$ws.Range("AA13").Value = -1
$ws.Range("AB13").Value = 1.1
$ws.Range("AC13").Value = 0.95
$ws.Range("AD13").Value = -1
# Row 17
$ws.Range("E17").Value = 'Incheon Utd'
# Row 19
$ws.Range("F19").Value = 'Gangwon FC'
# Row 24
$ws.Range("E24").Value = 'Gangwon FC'
# Row 25
$ws.Range("F25").Value = 'Incheon Utd'
# Row 31
$ws.Range("E31").Value = 'Gangwon FC'
$ws.Range("F31").Value = 'Incheon Utd'
# Row 33
$ws.Range("E33").Value = 'Incheon Utd'
# Row 35
$ws.Range("F35").Value = 'Gangwon FC'
# Row 38
$ws.Range("F38").Value = 'Incheon Utd'
# Row 40
$ws.Range("F40").Value = 'Gangwon FC'
# Row 45
$ws.Range("E45").Value = 'Incheon Utd'
# Row 49
$ws.Range("E49").Value = 'Gangwon FC'
# Row 53
$ws.Range("E53").Value = 'Incheon Utd'
# Row 54
$ws.Range("F54").Value = 'Gangwon FC'
# Row 56
$ws.Range("B56").Value = 6149879
$ws.Range("E56").Value = 'Gwangju FC'
$ws.Range("F56").Value = 'Daejeon Hana Citizen'
$ws.Range("G56").Value = 3
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 2
$ws.Range("K56").Value = 'H'
$ws.Range("L56").Value = 2.15
$ws.Range("M56").Value = 3.3
$ws.Range("N56").Value = 3.4
$ws.Range("O56").Value = 2.2
$ws.Range("P56").Value = 3
$ws.Range("Q56").Value = 3.6
$ws.Range("R56").Value = -0.25
$ws.Range("S56").Value = 1.9
$ws.Range("T56").Value = 1.95
$ws.Range("V56").Value = 2
$ws.Range("W56").Value = 1.85
$ws.Range("X56").Value = 1.2
$ws.Range("Y56").Value = -1
$ws.Range("AA56").Value = 0.8999999999999999
$ws.Range("AB56").Value = -1
$ws.Range("AC56").Value = 1
# Row 57
$ws.Range("B57").Value = 6149880
$ws.Range("E57").Value = 'FC Seoul'
$ws.Range("F57").Value = 'Pohang Steelers'
$ws.Range("G57").Value = 2
$ws.Range("H57").Value = 2
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 'D'
$ws.Range("L57").Value = 2.5
$ws.Range("M57").Value = 3.2
$ws.Range("N57").Value = 2.8
$ws.Range("O57").Value = 2.625
$ws.Range("P57").Value = 3.3
$ws.Range("Q57").Value = 2.6
$ws.Range("R57").Value = 0
$ws.Range("S57").Value = 1.95
$ws.Range("T57").Value = 1.9
$ws.Range("V57").Value = 1.825
$ws.Range("W57").Value = 2.025
$ws.Range("X57").Value = -1
$ws.Range("Y57").Value = 2.3
$ws.Range("AA57").Value = 0
$ws.Range("AB57").Value = 0
$ws.Range("AC57").Value = 0.825
# Row 60
$ws.Range("F60").Value = 'Gangwon FC'
# Row 61
$ws.Range("F61").Value = 'Incheon Utd'
# Row 63
$ws.Range("E63").Value = 'Incheon Utd'
# Row 66
$ws.Range("E66").Value = 'Gangwon FC'
# Row 68
$ws.Range("E68").Value = 'Gangwon FC'
# Row 70
$ws.Range("E70").Value = 'Incheon Utd'
# Row 74
$ws.Range("F74").Value = 'Gangwon FC'
# Row 77
$ws.Range("E77").Value = 'Incheon Utd'
# Row 81
$ws.Range("F81").Value = 'Incheon Utd'
# Row 83
$ws.Range("E83").Value = 'Gangwon FC'
# Row 86
$ws.Range("F86").Value = 'Incheon Utd'
# Row 87
$ws.Range("E87").Value = 'Gangwon FC'
# Row 94
$ws.Range("B94").Value = 6323586
$ws.Range("E94").Value = 'Suwon FC'
$ws.Range("F94").Value = 'Ulsan Hyundai'
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = 3
$ws.Range("J94").Value = 1
$ws.Range("L94").Value = 4.5
$ws.Range("M94").Value = 4
$ws.Range("N94").Value = 1.7
$ws.Range("O94").Value = 4.75
$ws.Range("P94").Value = 4
$ws.Range("Q94").Value = 1.666
$ws.Range("R94").Value = 0.75
$ws.Range("S94").Value = 2.025
$ws.Range("T94").Value = 1.825
$ws.Range("U94").Value = 3.25
$ws.Range("V94").Value = 2.05
$ws.Range("W94").Value = 1.8
$ws.Range("Z94").Value = 0.6659999999999999
$ws.Range("AA94").Value = -0.5
$ws.Range("AB94").Value = 0.4125
$ws.Range("AC94").Value = 1.05
$ws.Range("AD94").Value = -1
# Row 95
$ws.Range("B95").Value = 6323587
$ws.Range("E95").Value = 'Gwangju FC'
$ws.Range("F95").Value = 'Jeonbuk Motors'
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 1
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 2.7
$ws.Range("M95").Value = 3.4
$ws.Range("N95").Value = 2.45
$ws.Range("O95").Value = 2.6
$ws.Range("P95").Value = 3.3
$ws.Range("Q95").Value = 2.7
$ws.Range("R95").Value = 0
$ws.Range("S95").Value = 1.875
$ws.Range("T95").Value = 1.975
$ws.Range("U95").Value = 2.25
$ws.Range("V95").Value = 1.975
$ws.Range("W95").Value = 1.875
$ws.Range("Z95").Value = 1.7
$ws.Range("AA95").Value = -1
$ws.Range("AB95").Value = 0.9750000000000001
$ws.Range("AC95").Value = -1
$ws.Range("AD95").Value = 0.875
# Row 96
$ws.Range("E96").Value = 'Incheon Utd'
$ws.Range("F96").Value = 'Gangwon FC'
# Row 98
$ws.Range("B98").Value = 6353261
$ws.Range("E98").Value = 'Suwon FC'
$ws.Range("F98").Value = 'FC Seoul'
$ws.Range("G98").Value = 1
$ws.Range("H98").Value = 1
$ws.Range("L98").Value = 3.4
$ws.Range("M98").Value = 3.6
$ws.Range("N98").Value = 2
$ws.Range("O98").Value = 3.75
$ws.Range("P98").Value = 3.8
$ws.Range("Q98").Value = 1.85
$ws.Range("R98").Value = 0.5
$ws.Range("S98").Value = 2.025
$ws.Range("T98").Value = 1.825
$ws.Range("U98").Value = 2.75
$ws.Range("V98").Value = 1.825
$ws.Range("W98").Value = 2.025
$ws.Range("Y98").Value = 2.8
$ws.Range("AA98").Value = 1.025
$ws.Range("AB98").Value = -1
$ws.Range("AD98").Value = 1.025
# Row 99
$ws.Range("B99").Value = 6353260
$ws.Range("E99").Value = 'Pohang Steelers'
$ws.Range("F99").Value = 'Ulsan Hyundai'
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 0
$ws.Range("L99").Value = 2.375
$ws.Range("M99").Value = 3.5
$ws.Range("N99").Value = 2.8
$ws.Range("O99").Value = 2.55
$ws.Range("P99").Value = 3.4
$ws.Range("Q99").Value = 2.625
$ws.Range("R99").Value = 0
$ws.Range("S99").Value = 1.825
$ws.Range("T99").Value = 2.025
$ws.Range("U99").Value = 2.5
$ws.Range("V99").Value = 2.025
$ws.Range("W99").Value = 1.825
$ws.Range("Y99").Value = 2.4
$ws.Range("AA99").Value = 0
$ws.Range("AB99").Value = 0
$ws.Range("AD99").Value = 0.825
# Row 101
$ws.Range("E101").Value = 'Gangwon FC'
# Row 102
$ws.Range("E102").Value = 'Incheon Utd'
# Row 104
$ws.Range("B104").Value = 6387791
$ws.Range("E104").Value = 'FC Seoul'
$ws.Range("F104").Value = 'Jeonbuk Motors'
$ws.Range("G104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 'A'
$ws.Range("L104").Value = 2.375
$ws.Range("M104").Value = 3.3
$ws.Range("N104").Value = 2.875
$ws.Range("O104").Value = 2.4
$ws.Range("P104").Value = 3.3
$ws.Range("Q104").Value = 2.8
$ws.Range("R104").Value = 0
$ws.Range("S104").Value = 1.775
$ws.Range("T104").Value = 2.1
$ws.Range("U104").Value = 2.5
$ws.Range("V104").Value = 1.85
$ws.Range("W104").Value = 2
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = 1.8
$ws.Range("AB104").Value = 1.1
$ws.Range("AC104").Value = -1
$ws.Range("AD104").Value = 1
# Row 105
$ws.Range("B105").Value = 6384125
$ws.Range("E105").Value = 'Daegu FC'
$ws.Range("F105").Value = 'Suwon FC'
$ws.Range("G105").Value = 2
$ws.Range("I105").Value = 1
$ws.Range("J105").Value = 1
$ws.Range("K105").Value = 'D'
$ws.Range("L105").Value = 1.666
$ws.Range("M105").Value = 3.75
$ws.Range("N105").Value = 4.75
$ws.Range("O105").Value = 1.6
$ws.Range("P105").Value = 3.8
$ws.Range("Q105").Value = 5.25
$ws.Range("R105").Value = -0.75
$ws.Range("S105").Value = 1.85
$ws.Range("T105").Value = 2
$ws.Range("U105").Value = 2.75
$ws.Range("V105").Value = 2.025
$ws.Range("W105").Value = 1.825
$ws.Range("Y105").Value = 2.8
$ws.Range("Z105").Value = -1
$ws.Range("AB105").Value = 1
$ws.Range("AC105").Value = 1.025
$ws.Range("AD105").Value = -1
# Row 106
$ws.Range("B106").Value = 6384126
$ws.Range("E106").Value = 'Suwon Bluewings'
$ws.Range("F106").Value = 'Pohang Steelers'
$ws.Range("I106").Value = 1
$ws.Range("L106").Value = 4
$ws.Range("N106").Value = 1.909
$ws.Range("O106").Value = 4
$ws.Range("P106").Value = 3.3
$ws.Range("Q106").Value = 1.95
$ws.Range("R106").Value = 0.5
$ws.Range("S106").Value = 1.875
$ws.Range("T106").Value = 1.975
$ws.Range("U106").Value = 2.5
$ws.Range("V106").Value = 2.1
$ws.Range("W106").Value = 1.775
$ws.Range("X106").Value = 3
$ws.Range("AA106").Value = 0.875
$ws.Range("AD106").Value = 0.7749999999999999
# Row 107
$ws.Range("B107").Value = 6388625
$ws.Range("E107").Value = 'Ulsan Hyundai'
$ws.Range("F107").Value = 'Gangwon FC'
$ws.Range("G107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 'D'
$ws.Range("L107").Value = 1.75
$ws.Range("M107").Value = 3.75
$ws.Range("N107").Value = 4.333
$ws.Range("O107").Value = 2
$ws.Range("P107").Value = 3.4
$ws.Range("Q107").Value = 3.5
$ws.Range("R107").Value = -0.5
$ws.Range("S107").Value = 2.025
$ws.Range("T107").Value = 1.825
$ws.Range("V107").Value = 2
$ws.Range("W107").Value = 1.85
$ws.Range("X107").Value = -1
$ws.Range("Y107").Value = 2.4
$ws.Range("AA107").Value = -1
$ws.Range("AB107").Value = 0.825
$ws.Range("AD107").Value = 0.8500000000000001
# Row 108
$ws.Range("B108").Value = 6384418
$ws.Range("E108").Value = 'Gwangju FC'
$ws.Range("G108").Value = 1
$ws.Range("K108").Value = 'H'
$ws.Range("L108").Value = 1.909
$ws.Range("M108").Value = 3.4
$ws.Range("N108").Value = 4
$ws.Range("O108").Value = 1.95
$ws.Range("Q108").Value = 4
$ws.Range("S108").Value = 2
$ws.Range("T108").Value = 1.85
$ws.Range("U108").Value = 2.25
$ws.Range("V108").Value = 1.925
$ws.Range("W108").Value = 1.925
$ws.Range("X108").Value = 0.95
$ws.Range("Y108").Value = -1
$ws.Range("AA108").Value = 1
$ws.Range("AB108").Value = -1
$ws.Range("AD108").Value = 0.925
# Row 110
$ws.Range("F110").Value = 'Gangwon FC'
# Row 113
$ws.Range("F113").Value = 'Incheon Utd'
# Row 117
$ws.Range("F117").Value = 'Gangwon FC'
# Row 118
$ws.Range("E118").Value = 'Incheon Utd'
# Row 119
$ws.Range("B119").Value = 7334080
$ws.Range("E119").Value = 'Suwon Bluewings'
$ws.Range("F119").Value = 'Daejeon Hana Citizen'
$ws.Range("H119").Value = 2
$ws.Range("I119").Value = 2
$ws.Range("K119").Value = 'D'
$ws.Range("L119").Value = 2.625
$ws.Range("M119").Value = 3.4
$ws.Range("N119").Value = 2.55
$ws.Range("O119").Value = 2.6
$ws.Range("P119").Value = 3.4
$ws.Range("Q119").Value = 2.6
$ws.Range("R119").Value = 0
$ws.Range("S119").Value = 1.875
$ws.Range("T119").Value = 1.975
$ws.Range("U119").Value = 2.75
$ws.Range("V119").Value = 2
$ws.Range("W119").Value = 1.85
$ws.Range("X119").Value = -1
$ws.Range("Y119").Value = 2.4
$ws.Range("AA119").Value = 0
$ws.Range("AB119").Value = 0
$ws.Range("AC119").Value = 1
$ws.Range("AD119").Value = -1
# Row 120
$ws.Range("B120").Value = 7333491
$ws.Range("E120").Value = 'Ulsan Hyundai'
$ws.Range("F120").Value = 'Daegu FC'
$ws.Range("H120").Value = 0
$ws.Range("I120").Value = 0
$ws.Range("K120").Value = 'H'
$ws.Range("L120").Value = 1.909
$ws.Range("M120").Value = 3.5
$ws.Range("N120").Value = 4
$ws.Range("O120").Value = 1.727
$ws.Range("P120").Value = 3.75
$ws.Range("Q120").Value = 5
$ws.Range("R120").Value = -0.75
$ws.Range("S120").Value = 1.925
$ws.Range("T120").Value = 1.925
$ws.Range("U120").Value = 2.5
$ws.Range("V120").Value = 1.975
$ws.Range("W120").Value = 1.875
$ws.Range("X120").Value = 0.7270000000000001
$ws.Range("Y120").Value = -1
$ws.Range("AA120").Value = 0.925
$ws.Range("AB120").Value = -1
$ws.Range("AC120").Value = -1
$ws.Range("AD120").Value = 0.875
# Row 122
$ws.Range("F122").Value = 'Incheon Utd'
# Row 123
$ws.Range("B123").Value = 7334083
$ws.Range("E123").Value = 'Jeju United'
$ws.Range("F123").Value = 'FC Seoul'
$ws.Range("G123").Value = 0
$ws.Range("H123").Value = 0
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 2.7
$ws.Range("M123").Value = 3.3
$ws.Range("N123").Value = 2.55
$ws.Range("O123").Value = 2.75
$ws.Range("P123").Value = 3.3
$ws.Range("Q123").Value = 2.55
$ws.Range("R123").Value = 0
$ws.Range("S123").Value = 2.025
$ws.Range("T123").Value = 1.825
$ws.Range("U123").Value = 2.5
$ws.Range("V123").Value = 1.85
$ws.Range("W123").Value = 2
$ws.Range("Y123").Value = 2.3
$ws.Range("AA123").Value = 0
$ws.Range("AB123").Value = 0
$ws.Range("AC123").Value = -1
$ws.Range("AD123").Value = 1
# Row 124
$ws.Range("B124").Value = 7334075
$ws.Range("E124").Value = 'Daegu FC'
$ws.Range("F124").Value = 'Gwangju FC'
$ws.Range("G124").Value = 1
$ws.Range("H124").Value = 1
$ws.Range("I124").Value = 1
$ws.Range("J124").Value = 1
$ws.Range("L124").Value = 2.625
$ws.Range("M124").Value = 3.2
$ws.Range("N124").Value = 2.7
$ws.Range("O124").Value = 2.8
$ws.Range("P124").Value = 3.4
$ws.Range("Q124").Value = 2.45
$ws.Range("R124").Value = 0.25
$ws.Range("S124").Value = 1.8
$ws.Range("T124").Value = 2.05
$ws.Range("U124").Value = 2.25
$ws.Range("V124").Value = 2.05
$ws.Range("W124").Value = 1.8
$ws.Range("Y124").Value = 2.4
$ws.Range("AA124").Value = 0.4
$ws.Range("AB124").Value = -0.5
$ws.Range("AC124").Value = -0.5
$ws.Range("AD124").Value = 0.4
# Row 125
$ws.Range("B125").Value = 7333492
$ws.Range("E125").Value = 'Gangwon FC'
$ws.Range("F125").Value = 'Jeonbuk Motors'
$ws.Range("G125").Value = 1
$ws.Range("H125").Value = 1
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 'D'
$ws.Range("L125").Value = 2.75
$ws.Range("M125").Value = 3.2
$ws.Range("N125").Value = 2.6
$ws.Range("O125").Value = 3.2
$ws.Range("P125").Value = 3.25
$ws.Range("Q125").Value = 2.25
$ws.Range("R125").Value = 0.25
$ws.Range("S125").Value = 1.85
$ws.Range("T125").Value = 2
$ws.Range("U125").Value = 2.25
$ws.Range("V125").Value = 1.825
$ws.Range("W125").Value = 2.025
$ws.Range("Y125").Value = 2.25
$ws.Range("Z125").Value = -1
$ws.Range("AA125").Value = 0.425
$ws.Range("AB125").Value = -0.5
$ws.Range("AC125").Value = -0.5
$ws.Range("AD125").Value = 0.5125
# Row 126
$ws.Range("B126").Value = 7334084
$ws.Range("E126").Value = 'Suwon FC'
$ws.Range("F126").Value = 'Suwon Bluewings'
$ws.Range("G126").Value = 2
$ws.Range("H126").Value = 3
$ws.Range("I126").Value = 1
$ws.Range("J126").Value = 1
$ws.Range("K126").Value = 'A'
$ws.Range("L126").Value = 2.15
$ws.Range("M126").Value = 3.5
$ws.Range("N126").Value = 3.2
$ws.Range("O126").Value = 2.55
$ws.Range("P126").Value = 3.4
$ws.Range("Q126").Value = 2.625
$ws.Range("R126").Value = 0
$ws.Range("S126").Value = 1.875
$ws.Range("T126").Value = 1.975
$ws.Range("U126").Value = 2.75
$ws.Range("V126").Value = 1.9
$ws.Range("W126").Value = 1.95
$ws.Range("Y126").Value = -1
$ws.Range("Z126").Value = 1.625
$ws.Range("AA126").Value = -1
$ws.Range("AB126").Value = 0.9750000000000001
$ws.Range("AC126").Value = 0.8999999999999999
$ws.Range("AD126").Value = -1
# Row 128
$ws.Range("E128").Value = 'Gangwon FC'
# Row 131
$ws.Range("B131").Value = 7334086
$ws.Range("E131").Value = 'Incheon Utd'
$ws.Range("F131").Value = 'Suwon FC'
$ws.Range("G131").Value = 2
$ws.Range("I131").Value = 1
$ws.Range("L131").Value = 2.05
$ws.Range("M131").Value = 3.6
$ws.Range("N131").Value = 3.3
$ws.Range("O131").Value = 2.05
$ws.Range("R131").Value = -0.5
$ws.Range("S131").Value = 2.05
$ws.Range("T131").Value = 1.8
$ws.Range("V131").Value = 1.825
$ws.Range("W131").Value = 2.025
$ws.Range("X131").Value = 1.05
$ws.Range("AA131").Value = 1.05
$ws.Range("AC131").Value = -0.5
$ws.Range("AD131").Value = 0.5125
# Row 132
$ws.Range("B132").Value = 7334085
$ws.Range("E132").Value = 'FC Seoul'
$ws.Range("F132").Value = 'Suwon Bluewings'
$ws.Range("G132").Value = 0
$ws.Range("H132").Value = 1
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 'A'
$ws.Range("L132").Value = 1.65
$ws.Range("M132").Value = 4
$ws.Range("N132").Value = 5
$ws.Range("O132").Value = 1.727
$ws.Range("P132").Value = 3.8
$ws.Range("Q132").Value = 4.2
$ws.Range("R132").Value = -0.75
$ws.Range("S132").Value = 2
$ws.Range("T132").Value = 1.85
$ws.Range("U132").Value = 2.75
$ws.Range("V132").Value = 1.95
$ws.Range("W132").Value = 1.9
$ws.Range("X132").Value = -1
$ws.Range("Z132").Value = 3.2
$ws.Range("AA132").Value = -1
$ws.Range("AB132").Value = 0.8500000000000001
$ws.Range("AC132").Value = -1
$ws.Range("AD132").Value = 0.8999999999999999
# Row 133
$ws.Range("B133").Value = 7333496
$ws.Range("E133").Value = 'Pohang Steelers'
$ws.Range("F133").Value = 'Daegu FC'
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 0
$ws.Range("K133").Value = 'H'
$ws.Range("L133").Value = 1.85
$ws.Range("M133").Value = 3.4
$ws.Range("N133").Value = 3.6
$ws.Range("O133").Value = 2.1
$ws.Range("P133").Value = 3.3
$ws.Range("Q133").Value = 3.6
$ws.Range("R133").Value = -0.25
$ws.Range("S133").Value = 1.8
$ws.Range("T133").Value = 2.05
$ws.Range("U133").Value = 2.25
$ws.Range("V133").Value = 1.975
$ws.Range("W133").Value = 1.875
$ws.Range("X133").Value = 1.1
$ws.Range("Z133").Value = -1
$ws.Range("AA133").Value = 0.8
$ws.Range("AB133").Value = -1
$ws.Range("AD133").Value = 0.875
# Row 135
$ws.Range("F135").Value = 'Incheon Utd'
# Row 137
$ws.Range("F137").Value = 'Gangwon FC'
# Row 143
$ws.Range("E143").Value = 'Incheon Utd'
# Row 144
$ws.Range("E144").Value = 'Gangwon FC'
# Row 150
$ws.Range("F150").Value = 'Gangwon FC'
# Row 151
$ws.Range("F151").Value = 'Incheon Utd'
# Row 152
$ws.Range("F152").Value = 'Incheon Utd'
# Row 157
$ws.Range("F157").Value = 'Gangwon FC'
# Row 158
$ws.Range("E158").Value = 'Gangwon FC'
# Row 162
$ws.Range("E162").Value = 'Incheon Utd'
# Row 167
$ws.Range("E167").Value = 'Incheon Utd'
# Row 169
$ws.Range("F169").Value = 'Gangwon FC'
# Row 172
$ws.Range("E172").Value = 'Gangwon FC'
# Row 173
$ws.Range("B173").Value = 7716531
$ws.Range("E173").Value = 'Daegu FC'
$ws.Range("F173").Value = 'FC Seoul'
$ws.Range("G173").Value = 0
$ws.Range("H173").Value = 0
$ws.Range("I173").Value = 0
$ws.Range("J173").Value = 0
$ws.Range("K173").Value = 'D'
$ws.Range("L173").Value = 2.45
$ws.Range("M173").Value = 3.2
$ws.Range("N173").Value = 2.7
$ws.Range("O173").Value = 2.5
$ws.Range("P173").Value = 3.25
$ws.Range("Q173").Value = 2.8
$ws.Range("R173").Value = 0
$ws.Range("S173").Value = 1.85
$ws.Range("T173").Value = 2
$ws.Range("V173").Value = 1.9
$ws.Range("W173").Value = 1.95
$ws.Range("Y173").Value = 2.25
$ws.Range("Z173").Value = -1
$ws.Range("AA173").Value = 0
$ws.Range("AB173").Value = 0
$ws.Range("AC173").Value = -1
$ws.Range("AD173").Value = 0.95
# Row 174
$ws.Range("B174").Value = 7715278
$ws.Range("E174").Value = 'Jeonbuk Motors'
$ws.Range("F174").Value = 'Incheon Utd'
$ws.Range("G174").Value = 2
$ws.Range("H174").Value = 3
$ws.Range("I174").Value = 1
$ws.Range("J174").Value = 1
$ws.Range("K174").Value = 'A'
$ws.Range("L174").Value = 1.8
$ws.Range("M174").Value = 3.5
$ws.Range("N174").Value = 4.75
$ws.Range("O174").Value = 1.909
$ws.Range("P174").Value = 3.4
$ws.Range("Q174").Value = 4
$ws.Range("R174").Value = -0.5
$ws.Range("S174").Value = 1.975
$ws.Range("T174").Value = 1.875
$ws.Range("V174").Value = 1.95
$ws.Range("W174").Value = 1.9
$ws.Range("Y174").Value = -1
$ws.Range("Z174").Value = 3
$ws.Range("AA174").Value = -1
$ws.Range("AB174").Value = 0.875
$ws.Range("AC174").Value = 0.95
$ws.Range("AD174").Value = -1
# Row 179
$ws.Range("F179").Value = 'Incheon Utd'
# Row 181
$ws.Range("E181").Value = 'Gangwon FC'
# Row 185
$ws.Range("E185").Value = 'Incheon Utd'
$ws.Range("F185").Value = 'Gangwon FC'
# Row 188
$ws.Range("F188").Value = 'Incheon Utd'
# Row 190
$ws.Range("F190").Value = 'Gangwon FC'
# Row 197
$ws.Range("E197").Value = 'Incheon Utd'
# Row 198
$ws.Range("E198").Value = 'Gangwon FC'
# Row 201
$ws.Range("F201").Value = 'Incheon Utd'
# Row 202
$ws.Range("F202").Value = 'Gangwon FC'
# Row 205
$ws.Range("E205").Value = 'Gangwon FC'
# Row 206
$ws.Range("E206").Value = 'Incheon Utd'
# Row 209
$ws.Range("B209").Value = 7715299
$ws.Range("E209").Value = 'Ulsan Hyundai'
$ws.Range("F209").Value = 'Gimcheon Sangmu FC'
$ws.Range("G209").Value = 2
$ws.Range("H209").Value = 2
$ws.Range("J209").Value = 1
$ws.Range("L209").Value = 1.65
$ws.Range("M209").Value = 3.8
$ws.Range("N209").Value = 5.25
$ws.Range("Q209").Value = 4.5
$ws.Range("S209").Value = 1.95
$ws.Range("T209").Value = 1.9
$ws.Range("U209").Value = 2.75
$ws.Range("V209").Value = 1.875
$ws.Range("W209").Value = 1.975
$ws.Range("AB209").Value = 0.8999999999999999
$ws.Range("AC209").Value = 0.875
$ws.Range("AD209").Value = -1
# Row 210
$ws.Range("B210").Value = 7715301
$ws.Range("E210").Value = 'Pohang Steelers'
$ws.Range("F210").Value = 'Jeju United'
$ws.Range("G210").Value = 1
$ws.Range("H210").Value = 1
$ws.Range("J210").Value = 0
$ws.Range("L210").Value = 1.8
$ws.Range("M210").Value = 3.6
$ws.Range("N210").Value = 4.333
$ws.Range("Q210").Value = 4.75
$ws.Range("S210").Value = 1.925
$ws.Range("T210").Value = 1.925
$ws.Range("U210").Value = 2.5
$ws.Range("V210").Value = 1.8
$ws.Range("W210").Value = 2.05
$ws.Range("AB210").Value = 0.925
$ws.Range("AC210").Value = -1
$ws.Range("AD210").Value = 1.05
# Row 212
$ws.Range("F212").Value = 'Gangwon FC'
# Row 215
$ws.Range("B215").Value = 7715307
$ws.Range("E215").Value = 'Incheon Utd'
$ws.Range("F215").Value = 'Ulsan Hyundai'
$ws.Range("H215").Value = 0
$ws.Range("K215").Value = 'H'
$ws.Range("L215").Value = 3.3
$ws.Range("M215").Value = 3.4
$ws.Range("N215").Value = 2.15
$ws.Range("O215").Value = 3.6
$ws.Range("P215").Value = 3.8
$ws.Range("Q215").Value = 1.909
$ws.Range("R215").Value = 0.5
$ws.Range("S215").Value = 1.9
$ws.Range("T215").Value = 1.95
$ws.Range("U215").Value = 2.75
$ws.Range("V215").Value = 1.85
$ws.Range("W215").Value = 2
$ws.Range("X215").Value = 2.6
$ws.Range("Z215").Value = -1
$ws.Range("AA215").Value = 0.8999999999999999
$ws.Range("AB215").Value = -1
$ws.Range("AC215").Value = -1
$ws.Range("AD215").Value = 1
# Row 216
$ws.Range("B216").Value = 7716534
$ws.Range("E216").Value = 'FC Seoul'
$ws.Range("F216").Value = 'Daegu FC'
$ws.Range("H216").Value = 2
$ws.Range("K216").Value = 'A'
$ws.Range("L216").Value = 2.1
$ws.Range("M216").Value = 3.25
$ws.Range("N216").Value = 3.6
$ws.Range("O216").Value = 2.05
$ws.Range("P216").Value = 3.2
$ws.Range("Q216").Value = 3.9
$ws.Range("R216").Value = -0.5
$ws.Range("S216").Value = 2
$ws.Range("T216").Value = 1.85
$ws.Range("U216").Value = 2.25
$ws.Range("V216").Value = 1.95
$ws.Range("W216").Value = 1.9
$ws.Range("X216").Value = -1
$ws.Range("Z216").Value = 2.9
$ws.Range("AA216").Value = -1
$ws.Range("AB216").Value = 0.8500000000000001
$ws.Range("AC216").Value = 0.95
$ws.Range("AD216").Value = -1
# Row 221
$ws.Range("E221").Value = 'Gangwon FC'
# Row 222
$ws.Range("F222").Value = 'Incheon Utd'
# Row 226
$ws.Range("E226").Value = 'Incheon Utd'
# Row 229
$ws.Range("E229").Value = 'Gangwon FC'
# Row 231
$ws.Range("B231").Value = 7715318
$ws.Range("E231").Value = 'Gimcheon Sangmu FC'
$ws.Range("F231").Value = 'Pohang Steelers'
$ws.Range("L231").Value = 2.3
$ws.Range("M231").Value = 3.3
$ws.Range("N231").Value = 2.7
$ws.Range("O231").Value = 2.25
$ws.Range("P231").Value = 3.3
$ws.Range("Q231").Value = 2.75
$ws.Range("R231").Value = -0.25
$ws.Range("S231").Value = 2.025
$ws.Range("T231").Value = 1.825
$ws.Range("V231").Value = 2
$ws.Range("W231").Value = 1.85
$ws.Range("X231").Value = 1.25
$ws.Range("AA231").Value = 1.025
$ws.Range("AC231").Value = 1
# Row 232
$ws.Range("B232").Value = 7715317
$ws.Range("E232").Value = 'Suwon FC'
$ws.Range("F232").Value = 'Gangwon FC'
$ws.Range("L232").Value = 2.4
$ws.Range("M232").Value = 3.4
$ws.Range("N232").Value = 2.5
$ws.Range("O232").Value = 2.45
$ws.Range("P232").Value = 3.5
$ws.Range("Q232").Value = 2.4
$ws.Range("R232").Value = 0
$ws.Range("S232").Value = 1.975
$ws.Range("T232").Value = 1.875
$ws.Range("V232").Value = 1.975
$ws.Range("W232").Value = 1.875
$ws.Range("X232").Value = 1.45
$ws.Range("AA232").Value = 0.9750000000000001
$ws.Range("AC232").Value = 0.9750000000000001
# Row 235
$ws.Range("E235").Value = 'Incheon Utd'
# Row 236
$ws.Range("F236").Value = 'Gangwon FC'
